$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new price-report row was recorded for this market/variety and inserted
# at row 268, pushing the existing rows 268-283 down to 269-284.
$ws.Rows("268:268").Insert()

$ws.Range("A268").Value = 5
$ws.Range("B268").Value = "Macroferia Regional de Talca"
$ws.Range("C268").Value = "Maule"
$ws.Range("D268").Value = 44748
$ws.Range("E268").Value = 7
$ws.Range("F268").Value = "Fruta"
$ws.Range("G268").Value = 100102
$ws.Range("H268").Value = "Cítricos"
$ws.Range("I268").Value = 100102004
$ws.Range("J268").Value = "Mandarina"
$ws.Range("K268").Value = "Clemenuless"
$ws.Range("L268").Value = "Segunda"
$ws.Range("M268").Value = 200
$ws.Range("N268").Value = 5000
$ws.Range("O268").Value = 5000
$ws.Range("P268").Value = 5000
$ws.Range("Q268").Value = "`$/caja 18 kilos"
$ws.Range("R268").Value = "Región de O'Higgins"
$ws.Range("S268").Value = 278
$ws.Range("T268").Value = 18
